# CampusCrimesMerged.xlsx — add attribute-field header row to Sheet1
# (Sector of Institution, Reporting Location, Offense, Date, Count)
# and move the active selection to H10, matching the source commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Sector of Institution"
$ws.Range("B1").Value = "Reporting Location"
$ws.Range("C1").Value = "Offense"
$ws.Range("D1").Value = "Date"
$ws.Range("E1").Value = "Count"

$ws.Range("H10").Select() | Out-Null
